# Add new client row (row 2) to Sheet1, matching the existing header
# columns: ID | Création | Nom | Prénom | Mail | Date de naissance | Tél 1 | Tél 2 | Type | Notes
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Leading apostrophe forces these to be stored as literal text (not
# auto-converted to a number / date), so the ID, the French-formatted
# date and the phone number (with its leading zero) are preserved exactly
# as typed.
$ws.Range("A2").Value = "'11"
$ws.Range("B2").Value = "'21/08/2024"
$ws.Range("C2").Value = "Timothée"
$ws.Range("D2").Value = "Régis"
$ws.Range("E2").Value = "registimothee@gmail.com"
$ws.Range("F2").Value = "'08/08/1990"
$ws.Range("G2").Value = "'0769181771"
$ws.Range("H2").Value = ""
$ws.Range("I2").Value = "Nouveau client"
$ws.Range("J2").Value = ""
